{"js": "// The document's style sheet contains an unused \"Footnote Text\"\n// (styleId \"FootnoteText\") paragraph style that was dropped when the\n// docx was rebuilt. Remove it from the Styles collection.\nconst styles = context.document.getStyles();\nconst footnoteTextStyle = styles.getByNameOrNullObject(\"Footnote Text\");\nfootnoteTextStyle.load(\"nameLocal\");\nawait context.sync();\n\nif (!footnoteTextStyle.isNullObject) {\n  footnoteTextStyle.delete();\n  await context.sync();\n}\n", "ps1": "# The document's style sheet contains an unused \"Footnote Text\"\n# (styleId \"FootnoteText\") paragraph style that was dropped when the\n# docx was rebuilt. Remove it from the Styles collection.\n$d = $word.ActiveDocument\n\ntry {\n    $style = $d.Styles(\"Footnote Text\")\n    $style.Delete()\n} catch {\n    # Style already absent - nothing to do.\n}\n"}
